# Daily attendance processing - 2025-12-18 21:50:10
# Normalize "Recorded By" (column G) value ordering on the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact current value -> corrected (reordered) value.
$map = @{
    "dnasr281@gmail.com, System"            = "System, dnasr281@gmail.com";
    "System, admin@admin.com"               = "admin@admin.com, System";
    "backup@backdoor.com, System, system"   = "backup@backdoor.com, system, System";
    "dnasr281@gmail.com, admin@admin.com"   = "admin@admin.com, dnasr281@gmail.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
